# seed_purchases_purchaserequestheader.xlsx - update seed data
# - trims the purchase-request table down to the first two rows
# - refreshes their DATE / VENDOR / STATUS values
# - widens column A and moves the active selection
# - sets the sheet to portrait page orientation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-stale rows (P REQ 03/04/05) entirely - this also drops the
# shared strings that only they used and shrinks the sheet dimension.
$ws.Rows("4:6").Delete()

# Row 2 (P REQ 01): new date, vendor unchanged
$ws.Range("B2").Value = 45660

# Row 3 (P REQ 02): new date, new vendor
$ws.Range("B3").Value = 45660
$ws.Range("E3").Value = "Vendor 02"

# Status column refresh (order matters for shared-string table layout)
$ws.Range("F2").Value = "CLOSED"
$ws.Range("F3").Value = "CLOSED"

# Widen column A to fit the new content and move the selection
$ws.Columns("A").ColumnWidth = 7.5
[void]$ws.Range("F5").Select()

# Switch the sheet to portrait orientation
$ws.PageSetup.Orientation = 1
